# Add new application-log rows (39-48) to the "data" sheet, exercising
# the new date-extraction logic against several text layouts of the same
# delivery request (newline-separated, single typo date, no separators,
# single-space separated).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("R39").Value = '\nИГО:\nЗаявка на доставку \n1. Дата отгрузки 28.04.2023 \n2. Марка ЦЕМ I 42.5н Беларусь   \n3. Количество 40 тонн  \n4. От ООО Спарта \n5. Завод: Сзтк \n6. Покупатель ООО ""ТД"Цемент \n7. Грузополучатель: ООО "ТД"Цемент  \n8. Голицыно\n+7 910 404-06-14\nРБУ\nМожайское ш., 81\n'

$ws.Range("R40").Value = '\nИГО:\nЗаявка на доставку \n1. Дата отгрузки 28.04.2024 \n2. Марка ЦЕМ I 42.5н Беларусь   \n3. Количество 40 тонн  \n4. От ООО Спарта \n5. Завод: Сзтк \n6. Покупатель ООО ""ТД"Цемент \n7. Грузополучатель: ООО "ТД"Цемент  \n8. Голицыно\n+7 910 404-06-14\nРБУ\nМожайское ш., 81\n'

$ws.Range("R41").Value = '\nИГО:\nЗаявка на доставку \n1. Дата отгрузки 28.04.2024 \n2. Марка ЦЕМ I 42.5н Беларусь   \n3. Количество 40 тонн  \n4. От ООО Спарта \n5. Завод: Сзтк \n6. Покупатель ООО ""ТД"Цемент \n7. Грузополучатель: ООО "ТД"Цемент  \n8. Голицыно\n+7 910 404-06-14\nРБУ\nМожайское ш., 81\n'

$ws.Range("A42").Value = 'Игорь Хабаров'
$ws.Range("B42").Value = 'доставка'
$ws.Range("C42").Value = '28.04.2024'
$ws.Range("D42").Value = 'ЦЕМ I 42.5н Беларусь   '
$ws.Range("F42").Value = 40
$ws.Range("G42").Value = 'т'
$ws.Range("K42").Value = 'ООО ""ТД"Цемент '
$ws.Range("R42").Value = '\nИГО:\nЗаявка на доставку \n1. Дата отгрузки 28.04.2024 \n2. Марка ЦЕМ I 42.5н Беларусь   \n3. Количество 40 тонн  \n4. От ООО Спарта \n5. Завод: Сзтк \n6. Покупатель ООО ""ТД"Цемент \n7. Грузополучатель: ООО "ТД"Цемент  \n8. Голицыно\n+7 910 404-06-14\nРБУ\nМожайское ш., 81\n'

$ws.Range("A43").Value = 'Юрий'
$ws.Range("B43").Value = 'доставка'
$ws.Range("C43").Value = '26.04.2024'
$ws.Range("D43").Value = 'Щебень гранитный 5-20(ЛСР)'
$ws.Range("E43").Value = 'Обязательно актуальный паспорт!'
$ws.Range("F43").Value = 120
$ws.Range("G43").Value = 'т'
$ws.Range("K43").Value = 'ООО НВЛ ГРУП'
$ws.Range("R43").Value = '\nЮра Менеджер:\n1. Дата отгрузки:\n26.04.2024\n2. От ООО Спарта\n3. Марка: Щебень гранитный 5-20(ЛСР)\nОбязательно актуальный паспорт!\n4. Покупатель ООО НВЛ ГРУП\n5. Бетас \n6.  Количество 120т\n7. Машина: \nО327ВН790 МАN\nО039ОУ790 МАN,\nХ194ВА797 КАМАЗ\nВ683СН790 МАN\nМ991ХС750 МАN\nО030ХЕ123 МАN\nВ551ВО790 МАN\nС289УТ750 МАN\nВ247ХО750 МАN \nУ162ХК750 MAN\nХ240ВА797 КАМАЗ\n А215УХ750 МAN\n'

$ws.Range("A44").Value = 'Игорь Хабаров'
$ws.Range("B44").Value = 'доставка'
$ws.Range("C44").Value = '26.04.2024'
$ws.Range("D44").Value = 'ЦЕМ I 42.5Н БЦК '
$ws.Range("F44").Value = 35
$ws.Range("G44").Value = 'т'
$ws.Range("K44").Value = 'ООО "Бетонная индустрия»'
$ws.Range("R44").Value = '\nИГО:\n1. Дата отгрузки\n26.04.2024\n2.Марка цемента ЦЕМ I 42.5Н БЦК \n3. Количество 35 тонн \n4. Продажа от ООО "Спарта"\n5. С псо 13\n6. Покупатель ООО "Бетонная индустрия»\n7. Грузополучатель ООО "Бетонная индустрия»\n8. Адрес грузополучателя \nОдинцово \nКобяковская. Краснознаменск.\n'

$ws.Range("A45").Value = 'Игорь Хабаров'
$ws.Range("B45").Value = 'доставка'
$ws.Range("D45").Value = 'ЦЕМ I 42.5Н БЦК '
$ws.Range("F45").Value = 35
$ws.Range("G45").Value = 'т'
$ws.Range("K45").Value = 'ООО "Бетонная индустрия»'
$ws.Range("R45").Value = '\nИГО:\n1. Дата отгрузки\n26.4.2024\n2.Марка цемента ЦЕМ I 42.5Н БЦК \n3. Количество 35 тонн \n4. Продажа от ООО "Спарта"\n5. С псо 13\n6. Покупатель ООО "Бетонная индустрия»\n7. Грузополучатель ООО "Бетонная индустрия»\n8. Адрес грузополучателя \nОдинцово \nКобяковская. Краснознаменск.\n'

$ws.Range("A46").Value = 'Игорь Хабаров'
$ws.Range("B46").Value = 'доставка'
$ws.Range("C46").Value = '26.04.2024'
$ws.Range("D46").Value = 'ЦЕМ I 42.5Н БЦК '
$ws.Range("F46").Value = 35
$ws.Range("G46").Value = 'т'
$ws.Range("K46").Value = 'ООО "Бетонная индустрия»'
$ws.Range("R46").Value = '\nИГО:\n1. Дата отгрузки\n26.04.2024\n2.Марка цемента ЦЕМ I 42.5Н БЦК \n3. Количество 35 тонн \n4. Продажа от ООО "Спарта"\n5. С псо 13\n6. Покупатель ООО "Бетонная индустрия»\n7. Грузополучатель ООО "Бетонная индустрия»\n8. Адрес грузополучателя \nОдинцово \nКобяковская. Краснознаменск.\n'

$ws.Range("A47").Value = 'Игорь Хабаров'
$ws.Range("B47").Value = 'доставка'
$ws.Range("C47").Value = '26.04.2024'
$ws.Range("D47").Value = 'ЦЕМ I 42.5Н БЦК '
$ws.Range("F47").Value = 35
$ws.Range("G47").Value = 'т'
$ws.Range("K47").Value = 'ООО "Бетонная индустрия»'
$ws.Range("R47").Value = 'ИГО:1. Дата отгрузки26.04.20242.Марка цемента ЦЕМ I 42.5Н БЦК 3. Количество 35 тонн 4. Продажа от ООО "Спарта"5. С псо 136. Покупатель ООО "Бетонная индустрия»7. Грузополучатель ООО "Бетонная индустрия»8. Адрес грузополучателя Одинцово Кобяковская. Краснознаменск.'

$ws.Range("A48").Value = 'Игорь Хабаров'
$ws.Range("B48").Value = 'доставка'
$ws.Range("C48").Value = '26.04.2024'
$ws.Range("D48").Value = 'ЦЕМ I 42.5Н БЦК '
$ws.Range("F48").Value = 35
$ws.Range("G48").Value = 'т'
$ws.Range("K48").Value = 'ООО "Бетонная индустрия»'
$ws.Range("R48").Value = ' ИГО: 1. Дата отгрузки 26.04.2024 2.Марка цемента ЦЕМ I 42.5Н БЦК  3. Количество 35 тонн  4. Продажа от ООО "Спарта" 5. С псо 13 6. Покупатель ООО "Бетонная индустрия» 7. Грузополучатель ООО "Бетонная индустрия» 8. Адрес грузополучателя  Одинцово  Кобяковская. Краснознаменск. '
